# Solutions_and_states.xlsx — fill in the "Best solution (least steps)"
# column (C) for the first three result rows, and move the active
# selection to D6 (as left by the author when they saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 35
$ws.Range("C4").Value = 29
$ws.Range("C5").Value = 83

$ws.Range("D6").Select()
